# Applies the recorded edit: within each group of consecutive sighting rows
# (rows 2-4, 9-10, 18-20, 26-27), the per-sighting data (id, taxon info,
# locality, coordinates, accuracy, activity and public-comment fields) is
# cyclically rotated by one position: row N ends up holding the data that
# row N+1 held before the edit (the last row in a group wraps around to the
# first row's original data). Columns that are identical across every row in
# a group (dates/times, county/municipality, reporter, etc.) are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function ColLetterToIndex($letter) {
    $idx = 0
    foreach ($ch in $letter.ToCharArray()) {
        $idx = $idx * 26 + ([int][char]$ch - [int][char]'A' + 1)
    }
    return $idx
}

# Groups of rows whose "content" columns get cyclically rotated (row i <- row i+1,
# last row in the group <- first row in the group).
$groups = @(
    @(2, 3, 4),
    @(9, 10),
    @(18, 19, 20),
    @(26, 27)
)

# Only the columns that actually vary row-to-row within a sighting group need to
# move; everything else (dates, times, county, reporter, ...) is shared by the
# whole group and is left alone.
$cols = @("A", "B", "E", "F", "G", "H", "M", "P", "Q", "R", "S", "AC", "AE")
$colIdx = @{}
foreach ($c in $cols) { $colIdx[$c] = ColLetterToIndex $c }

foreach ($group in $groups) {
    $n = $group.Length

    # Snapshot "before" values for every relevant cell in this group first,
    # so writes to one row never affect what we read for another.
    $snapshot = @{}
    for ($i = 0; $i -lt $n; $i++) {
        $row = $group[$i]
        $rowVals = @{}
        foreach ($c in $cols) {
            $cell = $ws.Cells.Item($row, $colIdx[$c])
            $rowVals[$c] = $cell.Value2
        }
        $snapshot[$row] = $rowVals
    }

    # Write back: row at position i gets the snapshot from position i+1 (wrap).
    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $group[$i]
        $srcRow = $group[($i + 1) % $n]
        $srcVals = $snapshot[$srcRow]
        foreach ($c in $cols) {
            $val = $srcVals[$c]
            $cell = $ws.Cells.Item($destRow, $colIdx[$c])
            $isEmpty = ($null -eq $val) -or (($val -is [string]) -and ($val -eq ""))
            if ($isEmpty) {
                $cell.ClearContents()
            } else {
                $cell.Value = $val
            }
        }
    }
}
